# Insert a new weekly price-report row at row 158 on the only worksheet,
# shifting the existing rows 158..254 down to 159..255 (dimension grows to
# A1:R255). The new row captures a "Región del Maule" / "$/saco 20 kilos"
# observation dated 2022-01-21 (serial 44582).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 158 (and everything below it) down by one row.
$ws.Rows.Item(158).Insert()

# Populate the newly-inserted row 158 with the new record.
$ws.Cells.Item(158, 1).Value  = 10
$ws.Cells.Item(158, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(158, 3).Value  = "La Araucanía"
$ws.Cells.Item(158, 4).Value  = 44582
$ws.Cells.Item(158, 5).Value  = 9
$ws.Cells.Item(158, 6).Value  = 100114013
$ws.Cells.Item(158, 7).Value  = "Zanahoria"
$ws.Cells.Item(158, 8).Value  = "Sin especificar"
$ws.Cells.Item(158, 9).Value  = "Primera"
$ws.Cells.Item(158, 10).Value = 100
$ws.Cells.Item(158, 11).Value = 8000
$ws.Cells.Item(158, 12).Value = 8000
$ws.Cells.Item(158, 13).Value = 8000
$ws.Cells.Item(158, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(158, 15).Value = "Región del Maule"
$ws.Cells.Item(158, 16).Value = 400
$ws.Cells.Item(158, 17).Value = 20
$ws.Cells.Item(158, 18).Value = "Hortaliza"
